$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the rate text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.14 = 28807.5 pesos`n✅ 28807.5 pesos = 7.13 = 955.06 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update rate values ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 140
$ws2.Range("O10").Value = 4033.05
$ws2.Range("N12").Value = 4039
$ws2.Range("O12").Value = 133.906
